$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update postal code values in column A (rows 3-12) to reflect new set of post codes
$ws.Range("A3").Value = "EX17 3AH"
$ws.Range("A4").Value = "S70 1RU"
$ws.Range("A5").Value = "M1 1JQ"
$ws.Range("A6").Value = "L1 5AS"
$ws.Range("A7").Value = "M1 1AD"
$ws.Range("A8").Value = "HA8 7JL"
$ws.Range("A9").Value = "B17 9NP"
$ws.Range("A10").Value = "M1 4DZ"
$ws.Range("A11").Value = "TR18 2TN"
$ws.Range("A12").Value = "L1 3HD"

# Update the selected/active cell to A10
$ws.Range("A10").Select()

$wb.Save()
